# Reorder the comma-separated "Recorded By" names in column G.
# Every multi-name cell gets its name list reversed, except cells that
# already start with "dnasr281@gmail.com" (those are left untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # column G
    $value = $cell.Value()

    if ($value -and ($value -like "*,*")) {
        $parts = $value -split ", "

        if ($parts.Count -gt 1 -and $parts[0] -ne "dnasr281@gmail.com") {
            $reversed = $parts[($parts.Count - 1)..0]
            $newValue = [string]::Join(", ", $reversed)
            $cell.Value = $newValue
        }
    }
}
